{"js": "// Colors two \"tarefa\" (task) paragraphs to mark their new status, matching\n// the author's commit: the *Converter products mobile/normal... paragraph\n// becomes light-blue (00B0F0, the same color already used for other\n// \"in progress\" items in the list) and the *Converter o tostify em\n// responsivo paragraph becomes purple (7030A0).\n//\n// Paragraphs are located by searching for distinctive, stable substrings\n// instead of hard-coded indices so the script is resilient to unrelated\n// edits elsewhere in the document.\n\nconst body = context.document.body;\n\nconst productsSearch = body.search(\"*Converter products mobile e normal\", { matchCase: false });\nproductsSearch.load(\"items\");\nconst tostifySearch = body.search(\"*Converter o tostify em responsivo\", { matchCase: false });\ntostifySearch.load(\"items\");\n\nawait context.sync();\n\nif (productsSearch.items.length > 0) {\n  const productsParagraph = productsSearch.items[0].paragraphs.getFirst();\n  productsParagraph.font.color = \"#00B0F0\";\n}\n\nif (tostifySearch.items.length > 0) {\n  const tostifyParagraph = tostifySearch.items[0].paragraphs.getFirst();\n  tostifyParagraph.font.color = \"#7030A0\";\n}\n\nawait context.sync();\n", "ps1": "# Colors two \"tarefa\" (task) paragraphs to mark their new status, matching\n# the author's commit: the *Converter products mobile/normal... paragraph\n# becomes light-blue (00B0F0, the same color already used for other\n# \"in progress\" items in the list) and the *Converter o tostify em\n# responsivo paragraph becomes purple (7030A0).\n#\n# Paragraphs are located via Find (on distinctive, stable substrings)\n# instead of hard-coded indices so the script is resilient to unrelated\n# edits elsewhere in the document. Word's Font.Color is a BGR-packed long\n# (0x00BBGGRR), the reverse byte order of the OOXML <w:color w:val=\"RRGGBB\"/>\n# hex string, hence the 0xF0B000 / 0xA03070 literals below for RGB\n# 00B0F0 / 7030A0 respectively.\n\n$d = $word.ActiveDocument\n\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\"Converter products mobile e normal\")\nif ($found1) {\n    $productsParagraph = $rng1.Paragraphs(1)\n    $productsParagraph.Range.Font.Color = 0xF0B000\n}\n\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"Converter o tostify em responsivo\")\nif ($found2) {\n    $tostifyParagraph = $rng2.Paragraphs(1)\n    $tostifyParagraph.Range.Font.Color = 0xA03070\n}\n"}
